# Add the new "Alcoholic" class as row 11 of the "Game Clases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Alcoholic"
$ws.Range("C11").Value = "str"
$ws.Range("D11").Value = "dex"
$ws.Range("P11").Value = "Fighter"
$ws.Range("Q11").Value = "Blacksmith"
$ws.Range("R11").Value = 25
$ws.Range("S11").Value = 25
